$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add text labels in column F for rows 7, 8, 10, 11, 12, 13
# Shared string table must end up ordered: c-c nebo c-h, C-O, C=O, Al, Al2O3, Al(OH)3
$ws.Range("F7").Value = "c-c nebo c-h"
$ws.Range("F8").Value = "C-O"
$ws.Range("F10").Value = "C=O"
$ws.Range("F11").Value = "Al"
$ws.Range("F13").Value = "Al2O3"
$ws.Range("F12").Value = "Al(OH)3"

# Add numeric values in column E
$ws.Range("E11").Value = 72.8
$ws.Range("E13").Value = 74.5

# Add values/formulas in column G
$ws.Range("G11").Value = 72.8
$ws.Range("G12").Formula = "=B12-(B`$11-G`$11)"
$ws.Range("G13").Formula = "=B13-(B`$11-G`$11)"

# Update the active selection to D11
$ws.Range("D11").Select()
